# update test sheet to 1 file only, add TS delete customer and add script
# good method for filtering and verify delete
#
# 1) Rename the original (only) sheet "Sheet1" -> "List Name"
# 2) Add a "New Customer" sheet after it, with a small FirstName/LastName/
#    Postode table
# 3) Add a "Delete Customer" sheet after that, with a customerName list
#    used to test the delete-customer filtering/verify-delete script

$wb = $excel.ActiveWorkbook

# --- Sheet1 -> "List Name" -------------------------------------------------
$listName = $wb.Worksheets.Item(1)
$listName.Name = "List Name"

# Grab a copy of the existing body formatting (font/alignment) from the
# original sheet so the new sheets' cells pick up a matching look instead of
# the bare default style.
$listName.Range("A1").Copy()

# --- New Customer ------------------------------------------------------------
$newCustomer = $wb.Worksheets.Add($null, $listName)
$newCustomer.Name = "New Customer"

$newCustomer.Range("A1").Value = "FirstName"
$newCustomer.Range("B1").Value = "LastName"
$newCustomer.Range("C1").Value = "Postode"

$newCustomer.Range("A2").Value = "Putra"
$newCustomer.Range("B2").Value = "Alrasy"
$newCustomer.Range("C2").Value = 81200

$newCustomer.Range("A3").Value = "Lee"
$newCustomer.Range("B3").Value = "Chong"
$newCustomer.Range("C3").Value = 50000

$newCustomer.Range("A4").Value = "Ali"
$newCustomer.Range("B4").Value = "Muthu"
$newCustomer.Range("C4").Value = 64570

$newCustomerBody = $newCustomer.Range("A1:C4")
$newCustomerBody.PasteSpecial(-4122)
$newCustomerBody.Font.Name = "Arial"

$newCustomerPostcodes = $newCustomer.Range("C2:C4")
$newCustomerPostcodes.HorizontalAlignment = -4152

# --- Delete Customer -----------------------------------------------------
$deleteCustomer = $wb.Worksheets.Add($null, $newCustomer)
$deleteCustomer.Name = "Delete Customer"

$deleteCustomer.Range("A1").Value = "customerName"
$deleteCustomer.Range("A2").Value = "Harry"
$deleteCustomer.Range("A3").Value = "Ron"

$deleteCustomerBody = $deleteCustomer.Range("A1:A3")
$deleteCustomerBody.PasteSpecial(-4122)

# Leave the workbook selection the way it started: on the first sheet.
$listName.Activate()
